# Insert a new weekly price record above the current row 170.
# This shifts rows 170-204 down to 171-205 (preserving all their data and
# formatting), and the freshly inserted row 170 is populated with the new
# observation: same categorical/unit data as the (old) row 170, but an
# updated date (Fecha) and volume (Volumen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 170; existing rows 170.. shift down to 171..
$ws.Rows.Item(170).Insert()

# Seed the new row 170 with the data that used to live there (now at 171),
# then overwrite the two cells that actually change for the new record.
$ws.Range("A170:R170").Value2 = $ws.Range("A171:R171").Value2
$ws.Range("D170").Value2 = 45135
$ws.Range("J170").Value2 = 1100
